# rooms.xlsx — refresh the room/style lookup values, widen the new
# "style" helper column, and restore the sheet's scroll/selection +
# print orientation the way the author last left it in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Column B ("room"): the "Phòng làm việc" group (rows 17-22) was
#    renamed to "Phòng sinh hoạt chung".
# ---------------------------------------------------------------------
for ($row = 17; $row -le 22; $row++) {
    $ws.Cells.Item($row, 2).Value = "Phòng sinh hoạt chung"
}

# ---------------------------------------------------------------------
# 2) Column C ("style"): the style list was refreshed from the old
#    4-item rotation (Hiện đại, Scandinavian, Cổ điển, Tối giản) to a
#    new 5-item rotation, keyed off the row's id (column A).
# ---------------------------------------------------------------------
$styles = @("Cổ điển", "Bohemian", "Tối giản", "Bắc Âu", "Vintage")
for ($row = 2; $row -le 29; $row++) {
    $id = $row - 1
    $idx = ($id - 1) % 5
    $ws.Cells.Item($row, 3).Value = $styles[$idx]
}

# ---------------------------------------------------------------------
# 3) New helper column C needs an explicit width, matching the layout
#    once the style labels got longer (e.g. "Phòng sinh hoạt chung").
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 19.77734375

# ---------------------------------------------------------------------
# 4) Sheet view: scrolled down to row 19, with C32 now the active
#    selected cell (just past the last data row).
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$ws.Range("C32").Select()
$win.ScrollRow = 19
$win.ScrollColumn = 1

# ---------------------------------------------------------------------
# 5) Page setup: explicit portrait orientation for printing.
# ---------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
